$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "1"
$ws.Range("B6").Value = "Mayank Sharma"
$ws.Range("C6").Value = 1000
$ws.Range("D6").Value = "2025-09-14 21:27:56"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2"
$ws.Range("B7").Value = "Pooja Sharma"
$ws.Range("C7").Value = 1000
$ws.Range("D7").Value = "2025-09-14 21:29:04"
